$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 149, pushing existing rows 149-200 down to 150-201.
$ws.Rows.Item(149).Insert()

# Populate the new row 149 with the same constant columns as its neighbours,
# and the new data values for this record.
$ws.Cells.Item(149, 1).Value = 3
$ws.Cells.Item(149, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(149, 3).Value = "Coquimbo"
$ws.Cells.Item(149, 4).Value = 44468
$ws.Cells.Item(149, 4).NumberFormat = $ws.Cells.Item(150, 4).NumberFormat
$ws.Cells.Item(149, 5).Value = 5
$ws.Cells.Item(149, 6).Value = 100112009
$ws.Cells.Item(149, 7).Value = "Acelga"
$ws.Cells.Item(149, 8).Value = "Sin especificar"
$ws.Cells.Item(149, 9).Value = "Primera"
$ws.Cells.Item(149, 10).Value = 250
$ws.Cells.Item(149, 11).Value = 2000
$ws.Cells.Item(149, 12).Value = 2200
$ws.Cells.Item(149, 13).Value = 2096
$ws.Cells.Item(149, 14).Value = "$/docena de atados (6 kilos)"
$ws.Cells.Item(149, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(149, 16).Value = 349
$ws.Cells.Item(149, 17).Value = 6
$ws.Cells.Item(149, 18).Value = "Hortaliza"
